# Update regression-output cells on the active sheet to reflect the refreshed
# hourly-dataset model run (new Coef./Std.Err./t/P>|t|/CI/coef_pos figures).
$ws = $excel.ActiveWorkbook.ActiveSheet

$data = @{
    2 = @{ "B"=0; "H"=0.09640014858912448 }
    3 = @{ "B"=0.1042360274185795; "H"=0.2006361760077039 }
    4 = @{ "B"=0.1327415693877952; "H"=0.2291417179769197 }
    5 = @{ "B"=0.05557451227918586; "H"=0.1519746608683104 }
    6 = @{ "B"=0.04807634347780593; "C"=0.00483549525424892; "D"=6.608260337980272; "E"=0.01974750449272871; "F"=0.03857812139835367; "G"=0.05757456555725905; "H"=0.1444764920669304 }
    7 = @{ "B"=0.03666292457026844; "C"=0.004477224564120189; "D"=3.507778098800936; "E"=0.02329402308806026; "F"=0.0278539889973731; "G"=0.045471860143165; "H"=0.1330630731593929 }
    8 = @{ "B"=0.0324564869705281; "C"=0.003516804698453915; "D"=3.251608481702725; "E"=0.01987996194583503; "F"=0.02554780295339812; "G"=0.03936517098765764; "H"=0.1288566355596526 }
    9 = @{ "B"=0.03776673540075415; "C"=0.004015835282486879; "D"=3.783227656605955; "E"=0.01204012986149713; "F"=0.02987344789083541; "G"=0.0456600229106744; "H"=0.1341668839898786 }
    10 = @{ "B"=0.03063983820585928; "C"=0.00343710494745291; "D"=3.200323776621989; "E"=0.01224433255601132; "F"=0.02388544852430685; "G"=0.03739422788741251; "H"=0.1270399867949838 }
    11 = @{ "B"=0.03041860034738776; "H"=0.1268187489365122 }
    12 = @{ "B"=0.04662430579991968; "H"=0.1430244543890442 }
    13 = @{ "B"=0.05929880620547275; "H"=0.1556989547945972 }
    14 = @{ "B"=0.06313924993060988; "H"=0.1595393985197344 }
    15 = @{ "B"=0.06908009376956847; "H"=0.165480242358693 }
    16 = @{ "B"=0.07013845731300106; "H"=0.1665386059021255 }
    17 = @{ "B"=0.07314807562362598; "H"=0.1695482242127505 }
    18 = @{ "B"=-0.09640014858912448; "C"=0.0122303611188648; "D"=-17.20231439914257; "E"=0.06297925392800259; "F"=-0.1205017252643247; "G"=-0.0722985719139239; "H"=0 }
    19 = @{ "B"=0.07642607716927302; "H"=0.1728262257583975 }
    20 = @{ "B"=0.0812768117235523; "C"=0.01007387060976619; "D"=16.69061209008504; "E"=0.04909126925379184; "F"=0.06142199449032846; "G"=0.1011316289567763; "H"=0.1776769603126768 }
    21 = @{ "B"=0.08462383564843402; "H"=0.1810239842375585 }
    22 = @{ "B"=0.08721080934262318; "C"=0.01051495079948699; "D"=16.28591936377417; "E"=0.05111023187141769; "F"=0.06646040862272151; "G"=0.1079612100625253; "H"=0.1836109579317476 }
    23 = @{ "B"=0.09132598976443784; "C"=0.0110112893770334; "D"=3218145331650.303; "E"=0.05836741945072728; "F"=0.06958191521043852; "G"=0.1130700643184372; "H"=0.1877261383535623 }
    24 = @{ "B"=0.08908402062988305; "C"=0.01087746673390756; "D"=13.99310405047447; "E"=0.06595513263659947; "F"=0.06759855770788888; "G"=0.1105694835518773; "H"=0.1854841692190075 }
    25 = @{ "B"=0.08947289932199931; "C"=0.01060687481264844; "D"=12.76293819123774; "E"=0.070611966486368; "F"=0.06858118500825905; "G"=0.1103646136357392; "H"=0.1858730479111238 }
    26 = @{ "B"=0.09491631149826743; "C"=0.01110838118191226; "D"=12.84564315253356; "E"=0.07355619123896143; "F"=0.07301442864616031; "G"=0.1168181943503741; "H"=0.1913164600873919 }
    27 = @{ "B"=0.09465985696319772; "C"=0.01110201127043648; "D"=11.73945612444291; "E"=0.0758936972211408; "F"=0.07277628112562902; "G"=0.1165434328007666; "H"=0.1910600055523222 }
    28 = @{ "B"=0.1040658556852016; "C"=0.01128497505878369; "D"=11.61687582407435; "E"=0.1250278914080156; "F"=0.08183407517647785; "G"=0.1262976361939254; "H"=0.2004660042743261 }
    29 = @{ "B"=0.03635591221195728; "C"=0.003889794707221325; "D"=5.332568551475089; "E"=0.01785822230520033; "F"=0.02870167976068788; "G"=0.0440101446632265; "H"=0.1327560608010818 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
